$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'37.197.11"
$ws.Range("E2").Value = "  +1.47%  "

# Row 3
$ws.Range("D3").Value = "'2.021.93"
$ws.Range("E3").Value = "  +3.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'247.35"
$ws.Range("E5").Value = "  +1.24%  "

# Row 6
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "  +1.96%  "

# Row 7
$ws.Range("D7").Value = "'60.10"

# Row 9
$ws.Range("D9").Value = "'0.392"
$ws.Range("E9").Value = "  +4.36%  "

# Row 10
$ws.Range("D10").Value = "'0.0811"
$ws.Range("E10").Value = "  +2.34%  "

# Row 11
$ws.Range("E11").Value = "  +1.68%  "

# Row 12
$ws.Range("D12").Value = "'15.21"
$ws.Range("E12").Value = "  +6.59%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'22.43"
$ws.Range("E13").Value = "  +2.03%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.856"
$ws.Range("E14").Value = "  +2.86%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.317.51"
$ws.Range("E15").Value = "  +3.79%  "

# Row 17
$ws.Range("D17").Value = "'2.025.20"
$ws.Range("E17").Value = "  +3.41%  "

# Row 18
$ws.Range("D18").Value = "'37.097.15"
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").Value = "'70.62"
$ws.Range("E19").Value = "  +1.29%  "

# Row 20
$ws.Range("E20").Value = "  +1.71%  "

# Row 21
$ws.Range("D21").Value = "'5.26"
$ws.Range("E21").Value = "  +3.53%  "

# Row 22
$ws.Range("D22").Value = "'230.74"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$ws.Range("E25").Value = "  +0.70%  "

# Row 26
$ws.Range("D26").Value = "'9.47"
$ws.Range("E26").Value = "  +3.18%  "

# Row 27
$ws.Range("D27").Value = "'163.74"
$ws.Range("E27").Value = "  +2.02%  "

# Row 28
$ws.Range("D28").Value = "'0.137"
$ws.Range("E28").Value = "  -3.37%  "

# Row 29
$ws.Range("D29").Value = "'19.80"
$ws.Range("E29").Value = "  +2.05%  "

# Row 30
$ws.Range("E30").Value = "  +7.24%  "

# Row 31
$ws.Range("D31").Value = "'0.122"
$ws.Range("E31").Value = "  +2.25%  "

# Row 32
$ws.Range("D32").Value = "'4.83"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("D33").Value = "'0.0663"
$ws.Range("E33").Value = "  +7.72%  "

# Row 34
$ws.Range("D34").Value = "'4.53"
$ws.Range("E34").Value = "  +1.70%  "

# Row 35
$ws.Range("E35").Value = "  +8.20%  "

# Row 36
$ws.Range("D36").Value = "'3.47"
$ws.Range("E36").Value = "  -2.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.13%  "

# Row 38
$ws.Range("E38").Value = "  +1.87%  "

# Row 39
$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  -1.22%  "

# Row 40
$ws.Range("D40").Value = "'0.0983"
$ws.Range("E40").Value = "  +0.47%  "

# Row 41
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  +1.11%  "

# Row 42
$ws.Range("D42").Value = "'1.19"
$ws.Range("E42").Value = "  +1.72%  "

# Row 43
$ws.Range("D43").Value = "'0.0215"
$ws.Range("E43").Value = "  +1.95%  "

# Row 44
$ws.Range("D44").Value = "'16.72"
$ws.Range("E44").Value = "  +4.46%  "

# Row 45
$ws.Range("D45").Value = "'92.17"
$ws.Range("E45").Value = "  +4.04%  "

# Row 46
$ws.Range("D46").Value = "'1.393.15"
$ws.Range("E46").Value = "  +1.86%  "

# Row 47
$ws.Range("E47").Value = "  +2.94%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.20"
$ws.Range("E48").Value = "  +18.51%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.45"
$ws.Range("E49").Value = "  +4.51%  "

# Row 50
$ws.Range("E50").Value = "  +0.33%  "

# Row 51
$ws.Range("D51").Value = "'46.91"
$ws.Range("E51").Value = "  +3.28%  "
